# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2 through 11
$newValues = @{
    2  = 0
    3  = 1
    4  = 3
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 2
    10 = 3
    11 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
